$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @('ECs', 'Col11a1', 'Ddr1', 'ECs', 3, 1, 0.06886266666666667, 0.206588, 0.003217064976512115, 0.003217064976512115, 2, 0.6666666666666666, 0.166446, 0.4993379999999999, 0.008607068890887148, 0.00860706889088715, 0.011461915416, 0.103157238744, 0.00002768949987930001, 0.00002768949987930002),
    @('ECs', 'Col11a1', 'Ddr1', 'FAPs', 3, 1, 0.06886266666666667, 0.206588, 0.003217064976512115, 0.003217064976512115, 3, 1, 1.890286333333333, 5.670859, 0.09774836700492934, 0.09774836700492935, 0.1301701576768889, 1.171531419092, 0.0003144628480028105, 0.0003144628480028107),
    @('ECs', 'Col11a1', 'Ddr1', 'MuSCs', 3, 1, 0.06886266666666667, 0.206588, 0.003217064976512115, 0.003217064976512115, 3, 1, 17.13435266666667, 51.403058, 0.8860324297535294, 0.8860324297535294, 1.179917216233778, 10.619254946104, 0.00285042389781401, 0.00285042389781401),
    @('ECs', 'Col11a1', 'Ddr1', 'Resolving-Mac', 3, 1, 0.06886266666666667, 0.206588, 0.003217064976512115, 0.003217064976512115, 3, 1, 0.1472056666666667, 0.441617, 0.007612134350654087, 0.007612134350654088, 0.01013697475511111, 0.091232772796, 0.00002448873081599405, 0.00002448873081599406),
    @('FAPs', 'Col11a1', 'Ddr1', 'ECs', 3, 1, 19.329446, 57.988338, 0.903015911988821, 0.9030159119888211, 2, 0.6666666666666666, 0.166446, 0.4993379999999999, 0.008607068890887148, 0.00860706889088715, 3.217308968916, 28.955780720244, 0.007772320164055068, 0.007772320164055071),
    @('FAPs', 'Col11a1', 'Ddr1', 'FAPs', 3, 1, 19.329446, 57.988338, 0.903015911988821, 0.9030159119888211, 3, 1, 1.890286333333333, 5.670859, 0.09774836700492934, 0.09774836700492935, 36.53818760470467, 328.843688442342, 0.08826833077637425, 0.08826833077637428),
    @('FAPs', 'Col11a1', 'Ddr1', 'MuSCs', 3, 1, 19.329446, 57.988338, 0.903015911988821, 0.9030159119888211, 3, 1, 17.13435266666667, 51.403058, 0.8860324297535294, 0.8860324297535294, 331.1975446152894, 2980.777901537604, 0.8001013826055543, 0.8001013826055544),
    @('FAPs', 'Col11a1', 'Ddr1', 'Resolving-Mac', 3, 1, 19.329446, 57.988338, 0.903015911988821, 0.9030159119888211, 3, 1, 0.1472056666666667, 0.441617, 0.007612134350654087, 0.007612134350654088, 2.845403984727334, 25.608635862546, 0.006873878442837333, 0.006873878442837334),
    @('MuSCs', 'Col11a1', 'Ddr1', 'ECs', 3, 1, 1.981008666666667, 5.943026, 0.09254700563005056, 0.09254700563005058, 2, 0.6666666666666666, 0.166446, 0.4993379999999999, 0.008607068890887148, 0.00860706889088715, 0.329730968532, 2.967578716788, 0.0007965584531031658, 0.0007965584531031663),
    @('MuSCs', 'Col11a1', 'Ddr1', 'FAPs', 3, 1, 1.981008666666667, 5.943026, 0.09254700563005056, 0.09254700563005058, 3, 1, 1.890286333333333, 5.670859, 0.09774836700492934, 0.09774836700492935, 3.744673608814889, 33.702062479334, 0.009046318671533443, 0.009046318671533448),
    @('MuSCs', 'Col11a1', 'Ddr1', 'MuSCs', 3, 1, 1.981008666666667, 5.943026, 0.09254700563005056, 0.09254700563005058, 3, 1, 17.13435266666667, 51.403058, 0.8860324297535294, 0.8860324297535294, 33.94330113038978, 305.489710173508, 0.08199964826480725, 0.08199964826480728),
    @('MuSCs', 'Col11a1', 'Ddr1', 'Resolving-Mac', 3, 1, 1.981008666666667, 5.943026, 0.09254700563005056, 0.09254700563005058, 3, 1, 0.1472056666666667, 0.441617, 0.007612134350654087, 0.007612134350654088, 0.2916157014491111, 2.624541313042, 0.000704480240606685, 0.0007044802406066854),
    @('Resolving-Mac', 'Col11a1', 'Ddr1', 'ECs', 3, 1, 0.026115, 0.078345, 0.001220017404616152, 0.001220017404616152, 2, 0.6666666666666666, 0.166446, 0.4993379999999999, 0.008607068890887148, 0.00860706889088715, 0.00434673729, 0.03912063561, 0.00001050077384961256, 0.00001050077384961256),
    @('Resolving-Mac', 'Col11a1', 'Ddr1', 'FAPs', 3, 1, 0.026115, 0.078345, 0.001220017404616152, 0.001220017404616152, 3, 1, 1.890286333333333, 5.670859, 0.09774836700492934, 0.09774836700492935, 0.049364827595, 0.444283448355, 0.000119254709018821, 0.000119254709018821),
    @('Resolving-Mac', 'Col11a1', 'Ddr1', 'MuSCs', 3, 1, 0.026115, 0.078345, 0.001220017404616152, 0.001220017404616152, 3, 1, 17.13435266666667, 51.403058, 0.8860324297535294, 0.8860324297535294, 0.4474636198900001, 4.02717257901, 0.001080974985353644, 0.001080974985353644),
    @('Resolving-Mac', 'Col11a1', 'Ddr1', 'Resolving-Mac', 3, 1, 0.026115, 0.078345, 0.001220017404616152, 0.001220017404616152, 3, 1, 0.1472056666666667, 0.441617, 0.007612134350654087, 0.007612134350654088, 0.003844275985, 0.034598483865, 0.000009286936394074456, 0.000009286936394074458)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
